$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "date only" format C29 currently uses (YYYY-MM-DD) so the
# new row can reuse the exact same style further below.
$dateOnlyFormat = $ws.Range("C29").NumberFormat

# Row 29, column C currently uses the "date only" style (YYYY-MM-DD).
# The diff switches it to the "date+time" style (YYYY-MM-DD HH:MM:SS),
# matching the style used by the rows above it (e.g. C26/C27/C28).
$ws.Range("C29").NumberFormat = $ws.Range("C28").NumberFormat

# Append a brand-new row 30 with the next day's fuel-price data.
$ws.Range("A30").Value = 770.419
$ws.Range("B30").Value = 692.068
$ws.Range("C30").Value = 45758
# New row's date cell keeps the "date only" style that C29 used to have.
$ws.Range("C30").NumberFormat = $dateOnlyFormat
